$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the country label "Mt" -> "MT" for row 22 (Malta).
# A leading apostrophe preserves the cell's existing "quotePrefix" text style
# (the same xf index 4 used by all other region cells in column A).
$ws.Range("A22").Value = "'MT"

# Reflect the cell-selection change recorded in the saved file
$ws.Activate()
$ws.Range("A23").Select()
